$wb = $excel.ActiveWorkbook
$ws1 = $wb.ActiveSheet
$ws1.Name = "订阅的创意工坊"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "任务评价判定"

$ws2.Range("A1").Value = "任务评级"
$ws2.Range("B1").Value = "打成条件"

$ws2.Range("A2").Value = "完美"
$ws2.Range("B2").Value = "无货损及迟到"

$ws2.Range("A3").Value = "干得不错"
$ws2.Range("B3").Value = "货损小于2%的"

$ws2.Range("A4").Value = "姗姗来迟"
$ws2.Range("B4").Value = "迟到"

$ws2.Range("A5").Value = "差强人意"
$ws2.Range("B5").Value = "货损达2%未到5%的"

$ws2.Range("A6").Value = "合情合理"
$ws2.Range("B6").Value = "货损达5%未到10%的"

$ws2.Range("A7").Value = "尚可容忍"
$ws2.Range("B7").Value = "货损达10.1%未到19.9%的"

$ws2.Range("A8").Value = "残破不堪"
$ws2.Range("B8").Value = "货损达20%未到40%的"

$ws2.Range("A9").Value = "毛手毛脚"
$ws2.Range("B9").Value = "货损达40%极其以上的"

Write-Output "done"
